# Applies the commit's edits:
#  - "Correlation matrix": updated correlation values (rows 2-4, col G + row 3/4 recompute)
#  - "Toggles del" / "Toggles input del": collapsed from 4 data rows to 2 (B2=B3=128),
#    removing the old rows 4 & 5, and updated the two bar charts that plot those columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Correlation matrix sheet
# ---------------------------------------------------------------------------
$wsCorr = $wb.Worksheets.Item("Correlation matrix")

$wsCorr.Range("G2").Value = 0.06454972243679027

$wsCorr.Range("B3").Value = 0
$wsCorr.Range("D3").Value = 0
$wsCorr.Range("E3").Value = 0
$wsCorr.Range("F3").Value = 0
$wsCorr.Range("G3").Value = 0.06454972243679027

$wsCorr.Range("B4").Value = 0
$wsCorr.Range("D4").Value = 0
$wsCorr.Range("E4").Value = 0
$wsCorr.Range("F4").Value = 0
$wsCorr.Range("G4").Value = 0.06454972243679027

# ---------------------------------------------------------------------------
# 2) "Toggles del" sheet — update B2:B3, drop the old rows 4 & 5
# ---------------------------------------------------------------------------
$wsDel = $wb.Worksheets.Item("Toggles del")
$wsDel.Range("B2").Value = 128
$wsDel.Range("B3").Value = 128
$wsDel.Range("A4:B5").Delete()

$chartDel = $wsDel.ChartObjects().Item(1).Chart
$chartDel.SeriesCollection(1).Formula = "=SERIES(,,'Toggles del'!`$B`$2:`$B`$3,1)"

# ---------------------------------------------------------------------------
# 3) "Toggles input del" sheet — same shape of edit
# ---------------------------------------------------------------------------
$wsInputDel = $wb.Worksheets.Item("Toggles input del")
$wsInputDel.Range("B2").Value = 128
$wsInputDel.Range("B3").Value = 128
$wsInputDel.Range("A4:B5").Delete()

$chartInputDel = $wsInputDel.ChartObjects().Item(1).Chart
$chartInputDel.SeriesCollection(1).Formula = "=SERIES(,,'Toggles input del'!`$B`$2:`$B`$3,1)"
